$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 44838, 80, 22000, 22000, 22000, 1467),
    @(3, 44819, 70, 22000, 22000, 22000, 1467),
    @(4, 44365, 80, 25000, 25000, 25000, 1667),
    @(5, 44799, 80, 23000, 23000, 23000, 1533),
    @(6, 44782, 120, 24000, 24000, 24000, 1600),
    @(7, 44754, 90, 25000, 25000, 25000, 1667),
    @(8, 44792, 120, 24000, 24000, 24000, 1600),
    @(9, 44827, 90, 22000, 22000, 22000, 1467),
    @(10, 44775, 120, 24000, 24000, 24000, 1600),
    @(11, 44764, 90, 24000, 24000, 24000, 1600),
    @(12, 44831, 90, 25000, 25000, 25000, 1667),
    @(13, 44740, 90, 25000, 25000, 25000, 1667),
    @(14, 44750, 90, 25000, 25000, 25000, 1667),
    @(15, 44817, 90, 23000, 23000, 23000, 1533),
    @(16, 44810, 110, 22000, 22000, 22000, 1467),
    @(17, 44771, 90, 25000, 25000, 25000, 1667),
    @(18, 44778, 120, 24000, 24000, 24000, 1600),
    @(19, 44761, 100, 23000, 25000, 24000, 1600),
    @(20, 44806, 70, 23000, 23000, 23000, 1533),
    @(21, 44400, 80, 25000, 25000, 25000, 1667),
    @(22, 44757, 80, 25000, 25000, 25000, 1667),
    @(23, 44781, 70, 24000, 24000, 24000, 1600),
    @(24, 44789, 90, 24000, 24000, 24000, 1600),
    @(25, 44407, 90, 25000, 25000, 25000, 1667),
    @(26, 44418, 90, 25000, 25000, 25000, 1667),
    @(27, 44803, 90, 24000, 24000, 24000, 1600),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]   # D: Fecha
    $ws.Cells.Item($r, 10).Value = $row[2]  # J: Volumen
    $ws.Cells.Item($r, 11).Value = $row[3]  # K: Precio minimo
    $ws.Cells.Item($r, 12).Value = $row[4]  # L: Precio maximo
    $ws.Cells.Item($r, 13).Value = $row[5]  # M: Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $row[6]  # P: Precio $/Kg
}

Write-Output "applied"